$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.308.36"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4700"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06583"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.65"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08015"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.99"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.870.52"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.113"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6854"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.66"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.319.13"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007629"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.90%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.116.22"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.266"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.395"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.70"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.91"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.952"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09887"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.361"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.457"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04711"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7005"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.755"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.268"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.08"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.958"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4169"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8414"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.80"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.064"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.153"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "913.62"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.47"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05702"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.14%  "
